$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Daily Pivot" sheet: widen column G (7) and move the selection/view.
# ---------------------------------------------------------------------
$wsPivot = $wb.Worksheets.Item("Daily Pivot")
$wsPivot.Select()
$wsPivot.Columns.Item(7).ColumnWidth = 10.7
$excel.ActiveWindow.TopLeftCell = $wsPivot.Range("B2")
$wsPivot.Range("J22").Select()

# ---------------------------------------------------------------------
# "Daily Expenditure" sheet: the bulk of the new data.
# ---------------------------------------------------------------------
$wsExp = $wb.Worksheets.Item("Daily Expenditure")
$wsExp.Select()

# Row 4 - debt partially repaid: record payment (T/U) and payoff date (V),
# and update the running-balance formula to subtract both T and U.
$wsExp.Range("K4").Copy()
$wsExp.Range("V4").PasteSpecial(-4122)
$wsExp.Range("T4").Value = 24000
$wsExp.Range("U4").Value = 63500
$wsExp.Range("V4").Value = 44667
$wsExp.Range("W4").Formula = "=N4+S4-T4-U4"

# Row 5 - same treatment as row 4.
$wsExp.Range("K5").Copy()
$wsExp.Range("V5").PasteSpecial(-4122)
$wsExp.Range("T5").Value = 24000
$wsExp.Range("U5").Value = 63500
$wsExp.Range("V5").Value = 44667
$wsExp.Range("W5").Formula = "=N5+S5-T5-U5"

# Row 14 - give the balance cell an explicit formula (value stays 0).
$wsExp.Range("W14").Formula = "=N14+T14-U14"

# Row 47 - Paddy sales income revised upward.
$wsExp.Range("H47").Value = 51500

# Row 51 - new expense line: Areca Plant purchase for 5000.
$wsExp.Range("F51").Value = "Areca"
$wsExp.Range("G51").Value = "Areca Plant"
$wsExp.Range("H51").Value = 5000

# Rows 54 & 55 - two new blank "Expense" rows appended, copied from 52:53.
$wsExp.Range("A52:H53").Copy($wsExp.Range("A54:H55"))
$wsExp.Range("A54").Value = 54
$wsExp.Range("A55").Value = 55
$wsExp.Range("C54").Formula = '=TEXT(B54, "mmm")'
$wsExp.Range("D54").Formula = '=TEXT(B54, "yyyy")'
$wsExp.Range("C55").Formula = '=TEXT(B55, "mmm")'
$wsExp.Range("D55").Formula = '=TEXT(B55, "yyyy")'

$wsExp.Range("B57").Select()
